{"js": "// Insert the new \"map dashboards\" draft notes after the existing\n// \"Are tracts with a high percentage...\" paragraph, right before the\n// end of the document body (sectPr).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Anchor on the last paragraph currently in the document body.\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\nconst newLines = [\n  \"\",\n  \"Here are income levels by census tract\",\n  \"Side by side maps\",\n  \"Joining things at county level\",\n  \"\",\n  \"I can maybe map with the tracts if I find geospatial data\",\n  \"County map could be helpful \\u2013 county polygons\",\n];\n\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Insert the new \"map dashboards\" draft notes after the existing\n# \"Are tracts with a high percentage...\" paragraph, at the end of the\n# document body (before sectPr).\n$d = $word.ActiveDocument\n\n$lines = @(\n    \"\",\n    \"Here are income levels by census tract\",\n    \"Side by side maps\",\n    \"Joining things at county level\",\n    \"\",\n    \"I can maybe map with the tracts if I find geospatial data\",\n    (\"County map could be helpful \" + [char]0x2013 + \" county polygons\")\n)\n\n$r = $d.Content\n$r.Collapse(0)  # wdCollapseEnd -> end of document content\n\nforeach ($line in $lines) {\n    $r.InsertAfter(\"`r\" + $line)\n    $r.Collapse(0)\n}\n"}
